$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the 24-pin female header part with a better one (less tight fit) ---
$ws.Range("F32").Value = "PPTC241LFBN-RC"
$ws.Range("G32").Value = "S7022-ND"

# New unit cost for the replacement part (I32 / I37 totals recalc automatically)
$ws.Range("H32").Value = 1.37

# --- Rebuild hyperlinks so the G32 Digikey-part-number hyperlink shows the new
#     display text while every other hyperlink (target + display) is preserved
#     exactly as it was. (Mutating a single Hyperlink object in place isn't
#     supported by this runtime - it silently duplicates the link instead of
#     replacing it - so the whole collection is recreated from scratch.) ---
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G2"), "https://www.digikey.com/product-detail/en/panasonic-electronic-components/EVQ-PF303M/P12212S-ND/593395", "", "", "P12212S-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.digikey.com/product-detail/en/panasonic-electronic-components/ERA-3ARW302V/P3.0KBECT-ND/3073316", "", "", "P3.0KBECT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.digikey.com/product-detail/en/tdk-lambda-americas-inc/CC3-0512DF-E/445-2465-ND/920425", "", "", "445-2465-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G5"), "https://www.digikey.com/product-detail/en/samsung-electro-mechanics/CL10B105KO8NNNC/1276-1019-1-ND/3889105", "", "", "1276-1019-1-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G6"), "https://www.digikey.com/product-detail/en/on-semiconductor/MC78L05ACDX/MC78L05ACDXCT-ND/3042647", "", "", "MC78L05ACDXCT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G7"), "https://www.digikey.com/product-detail/en/te-connectivity-amp-connectors/1-2834016-2/A123829-ND/5872945", "", "", "A123829-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G8"), "https://www.digikey.com/product-detail/en/samsung-electro-mechanics/CL10F104ZO8NNNC/1276-1258-1-ND/3889344", "", "", "CL10F104ZO8NNNC") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G9"), "https://www.digikey.com/product-detail/en/wurth-electronics-inc/885012206070/732-7988-1-ND/5454615", "", "", "732-7988-1-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G10"), "https://www.digikey.com/product-detail/en/yageo/CC1206ZKY5V7BB106/311-1376-1-ND/2103160", "", "", "311-1376-1-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G11"), "https://www.digikey.com/product-detail/en/texas-instruments/ISO7760DWR/296-48142-1-ND/8347457", "", "", "296-48142-1-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G12"), "https://www.digikey.com/product-detail/en/texas-instruments/ISO7341FCDWR/296-47779-1-ND/8133128", "", "", "296-47779-1-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G13"), "https://www.digikey.com/product-detail/en/analog-devices-inc/AD7321BRUZ-REEL7/AD7321BRUZ-REEL7CT-ND/4909360", "", "", "AD7321BRUZ-REEL7CT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G14"), "https://www.digikey.com/product-detail/en/analog-devices-inc/AD5752AREZ/AD5752AREZ-ND/1979362", "", "", "AD5752AREZ-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G15"), "https://www.digikey.com/product-detail/en/analog-devices-inc/AD1582ARTZ-REEL7/AD1582ARTZREEL7CT-ND/751228", "", "", "AD1582ARTZREEL7CT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G16"), "https://www.digikey.com/product-detail/en/texas-instruments/OPA197IDR/296-43866-1-ND/5880541", "", "", "296-43866-1-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G17"), "https://www.digikey.com/product-detail/en/visual-communications-company-vcc/5300H5/L20015-ND/59969", "", "", "L20015-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G18"), "https://www.digikey.com/product-detail/en/on-shore-technology-inc/302-R161/ED10535-ND/2794246", "", "", "ED10535-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G19"), "https://www.digikey.com/product-detail/en/sparkfun-electronics/DEV-14055/1568-1443-ND/6235191", "", "", "1568-1443-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G20"), "https://www.digikey.com/product-detail/en/molex-llc/0731375003/WM5514-ND/1465136", "", "", "WM5514-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G21"), "https://www.digikey.com/product-detail/en/susumu/RR0816P-102-D/RR08P1.0KDCT-ND/432724", "", "", "RR08P1.0KDCT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G22"), "https://www.digikey.com/product-detail/en/susumu/RR0816P-103-D/RR08P10.0KDCT-ND/432748", "", "", "RR08P10.0KDCT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G23"), "https://www.digikey.com/product-detail/en/yageo/RT0603DRE075K05L/311-2629-1-ND/6129048", "", "", "311-2629-1-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G24"), "https://www.digikey.com/product-detail/en/bourns-inc/CHV1206-JW-105ELF/CHV1206-JW-105ELFCT-ND/5176007", "", "", "CHV1206-JW-105ELFCT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G25"), "https://www.digikey.com/product-detail/en/kemet/C1206C222MGRAC7800/399-13198-1-ND/5879352", "", "", "399-13198-1-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G26"), "https://www.digikey.com/product-detail/en/susumu/RR0816P-182-D/RR08P1.8KDCT-ND/432730", "", "", "RR08P1.8KDCT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G27"), "https://www.digikey.com/product-detail/en/susumu/RR0816P-101-D/RR08P100DCT-ND/432700", "", "", "RR08P100DCT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G28"), "https://www.digikey.com/product-detail/en/analog-devices-inc/AD8276ARZ/AD8276ARZ-ND/2057775", "", "", "AD8276ARZ-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G29"), "https://www.digikey.com/product-detail/en/texas-instruments/INA826AIDR/296-30238-1-ND/3045446", "", "", "296-30238-1-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G30"), "https://www.digikey.com/product-detail/en/vishay-siliconix/DG509BEY-T1-E3/DG509BEY-T1-E3CT-ND/2296899", "", "", "DG509BEY-T1-E3CT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J30"), "https://www.digikey.com/product-detail/en/vishay-siliconix/DG409DY-T1-E3/DG409DY-T1-E3CT-ND/1850070", "", "", "Alternate part: DG409DY-T1-E3CT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G31"), "https://www.digikey.com/product-detail/en/susumu/RR0816P-201-D/RR08P200DCT-ND/432707", "", "", "RR08P200DCT-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G32"), "https://www.digikey.com/product-detail/en/samtec-inc/SSQ-124-03-T-S/SAM1206-24-ND/1111934", "", "", "S7022-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G33"), "https://www.digikey.com/product-detail/en/sullins-connector-solutions/PREC024SAAN-RC/S1012EC-24-ND/2774830", "", "", "S1012EC-24-ND ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G34"), "https://www.digikey.com/product-detail/en/bud-industries/EXN-23359-SVP/377-2550-ND/5886377", "", "", "377-2550-ND") | Out-Null

# Re-adding the hyperlinks resets the touched cells to Excel's built-in blue
# underlined "Hyperlink" font; restore each cell's original (non-underlined)
# font so only the intended G32 content actually changes visually.
$fontSpecs = @{
    "G2" = @("Arial", 16711680)
    "G3" = @("Times New Roman", 16711680)
    "G4" = @("Times New Roman", 16711680)
    "G5" = @("Times New Roman", 16711680)
    "G6" = @("Times New Roman", 16711680)
    "G7" = @("Times New Roman", 16711680)
    "G8" = @("Times New Roman", 16711680)
    "G9" = @("Times New Roman", 16711680)
    "G10" = @("Times New Roman", 16711680)
    "G11" = @("Times New Roman", 16711680)
    "G12" = @("Times New Roman", 16711680)
    "G13" = @("Times New Roman", 16711680)
    "G14" = @("Times New Roman", 16711680)
    "G15" = @("Times New Roman", 16711680)
    "G16" = @("Times New Roman", 16711680)
    "G17" = @("Arial", 16711680)
    "G18" = @("Times New Roman", 16711680)
    "G19" = @("Arial", 16711680)
    "G20" = @("Times New Roman", 16711680)
    "G21" = @("Times New Roman", 16711680)
    "G22" = @("Arial", 16711680)
    "G23" = @("Arial", 16711680)
    "G24" = @("Arial", 16711680)
    "G25" = @("Arial", 16711680)
    "G26" = @("Times New Roman", 16711680)
    "G27" = @("Arial", 16711680)
    "G28" = @("Times New Roman", 16711680)
    "G29" = @("Times New Roman", 16711680)
    "G30" = @("Times New Roman", 16711680)
    "J30" = @("Arial", 0)
    "G31" = @("Arial", 16711680)
    "G32" = @("Arial", 16711680)
    "G33" = @("Arial", 16711680)
    "G34" = @("Arial", 0)
}
foreach ($addr in $fontSpecs.Keys) {
    $spec = $fontSpecs[$addr]
    $c = $ws.Range($addr)
    $c.Font.Name = $spec[0]
    $c.Font.Size = 10
    $c.Font.Color = $spec[1]
    $c.Font.Underline = 0
}

# --- Update the view: scroll position and active cell selection ---
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H25").Select()
